$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Update existing rows 13-22 (weekly price-sheet refresh: values shift
# down a week, newest week's figures land in row 13).
# Only the cells that actually changed per row are touched.
# ---------------------------------------------------------------------

# Row 13
$ws.Cells.Item(13, 4).Value = 44427
$ws.Cells.Item(13, 10).Value = 360

# Row 14
$ws.Cells.Item(14, 4).Value = 44413
$ws.Cells.Item(14, 10).Value = 1200
$ws.Cells.Item(14, 11).Value = 10000
$ws.Cells.Item(14, 12).Value = 11000
$ws.Cells.Item(14, 13).Value = 10500
$ws.Cells.Item(14, 16).Value = 420

# Row 15
$ws.Cells.Item(15, 4).Value = 44377
$ws.Cells.Item(15, 10).Value = 800
$ws.Cells.Item(15, 11).Value = 9000
$ws.Cells.Item(15, 12).Value = 10000
$ws.Cells.Item(15, 13).Value = 9500
$ws.Cells.Item(15, 16).Value = 380

# Row 16
$ws.Cells.Item(16, 4).Value = 44426
$ws.Cells.Item(16, 10).Value = 500
$ws.Cells.Item(16, 11).Value = 11000
$ws.Cells.Item(16, 12).Value = 12000
$ws.Cells.Item(16, 13).Value = 11500
$ws.Cells.Item(16, 16).Value = 460

# Row 17
$ws.Cells.Item(17, 4).Value = 44412
$ws.Cells.Item(17, 10).Value = 1000
$ws.Cells.Item(17, 12).Value = 11000
$ws.Cells.Item(17, 13).Value = 10500
$ws.Cells.Item(17, 16).Value = 420

# Row 18
$ws.Cells.Item(18, 4).Value = 44336
$ws.Cells.Item(18, 10).Value = 1200
$ws.Cells.Item(18, 11).Value = 12000
$ws.Cells.Item(18, 12).Value = 13000
$ws.Cells.Item(18, 13).Value = 12500
$ws.Cells.Item(18, 16).Value = 500

# Row 19
$ws.Cells.Item(19, 4).Value = 44371
$ws.Cells.Item(19, 11).Value = 10000
$ws.Cells.Item(19, 12).Value = 12000
$ws.Cells.Item(19, 13).Value = 11000
$ws.Cells.Item(19, 16).Value = 440

# Row 20
$ws.Cells.Item(20, 4).Value = 44364
$ws.Cells.Item(20, 10).Value = 700
$ws.Cells.Item(20, 11).Value = 11000
$ws.Cells.Item(20, 12).Value = 12000
$ws.Cells.Item(20, 13).Value = 11500
$ws.Cells.Item(20, 16).Value = 460

# Row 21
$ws.Cells.Item(21, 4).Value = 44391
$ws.Cells.Item(21, 10).Value = 500

# Row 22
$ws.Cells.Item(22, 4).Value = 44399

# ---------------------------------------------------------------------
# Append four new weekly records as rows 23-26.
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row = 23; Date = 44419; Vol = 1100; Min = 11000; Max = 12000; Avg = 11500; PKg = 460 },
    @{ Row = 24; Date = 44420; Vol = 1000; Min = 10000; Max = 11000; Avg = 10500; PKg = 420 },
    @{ Row = 25; Date = 44398; Vol = 400;  Min = 9000;  Max = 10000; Avg = 9500;  PKg = 380 },
    @{ Row = 26; Date = 44343; Vol = 500;  Min = 9000;  Max = 10000; Avg = 9500;  PKg = 380 }
)

foreach ($rec in $newRows) {
    $r = $rec.Row
    $ws.Cells.Item($r, 1).Value = 2
    $ws.Cells.Item($r, 2).Value = 'Comercializadora del Agro de Limarí'
    $ws.Cells.Item($r, 3).Value = 'Coquimbo'
    $ws.Cells.Item($r, 4).Value = $rec.Date
    $ws.Cells.Item($r, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = 100112026
    $ws.Cells.Item($r, 7).Value = 'Haba'
    $ws.Cells.Item($r, 8).Value = 'Sin especificar'
    $ws.Cells.Item($r, 9).Value = 'Primera'
    $ws.Cells.Item($r, 10).Value = $rec.Vol
    $ws.Cells.Item($r, 11).Value = $rec.Min
    $ws.Cells.Item($r, 12).Value = $rec.Max
    $ws.Cells.Item($r, 13).Value = $rec.Avg
    $ws.Cells.Item($r, 14).Value = '$/saco 25 kilos'
    $ws.Cells.Item($r, 15).Value = 'Provincia de Limarí'
    $ws.Cells.Item($r, 16).Value = $rec.PKg
    $ws.Cells.Item($r, 17).Value = 25
    $ws.Cells.Item($r, 18).Value = 'Hortaliza'
}
